$p = $ppt.ActivePresentation

foreach ($s in $p.Slides) {
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Name -eq "Straight Connector 29") {
            $shp.Delete()
        }
    }
}
